# Fixed date error on resume
#
# The "Developer Intern" line listed the wrong start year (2017 instead
# of 2018). Fixing the year means placing the caret in the middle of
# that sentence and retyping the digit, which is exactly why Word ends
# up splitting that run into three pieces and re-anchoring its "last
# edit" bookmark (_GoBack) at the point where the user stopped typing,
# right after the corrected year.

$d = $word.ActiveDocument

function Get-ParagraphRangeContaining($doc, $needle) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p.Range
        }
    }
    return $null
}

# --- 1. Drop the old _GoBack bookmark. Word keeps only one at a time,
#        and it is about to be re-created at the real edit site below. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Correct the year, scoped to the Developer Intern paragraph so
#        the other (already-correct) "JAN 2017"-less dates elsewhere in
#        the resume are left untouched. ---
$introRange = Get-ParagraphRangeContaining $d "Developer Intern"
$introRange.Find.Execute("JAN 2017", $true, $false, $false, $false, $false, $true, 1, $false, "JAN 2018", 2)

# --- 3. Re-create the run split that happens where the caret first
#        landed ("Entropy Mu|ltimedia..."). A zero-length bookmark
#        forces Word to break the run at that character boundary; we
#        immediately discard the temporary bookmark since it isn't part
#        of the final document. ---
$introRange = Get-ParagraphRangeContaining $d "Developer Intern"
$introRange.Find.Execute("Entropy Mu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $introRange.End
$d.Bookmarks.Add("ZZZTempRunSplit", $d.Range($splitPoint, $splitPoint))
$d.Bookmarks("ZZZTempRunSplit").Delete()

# --- 4. Put _GoBack where the user actually finished typing: right
#        after the corrected year, before " - present". ---
$introRange = Get-ParagraphRangeContaining $d "Developer Intern"
$introRange.Find.Execute("JAN 2018", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPoint = $introRange.End
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPoint, $goBackPoint))

# --- 5. The same proofing pass flagged "Abundez" (part of
#        "Abundez-Arce") as an unrecognized word, which splits that run
#        at the hyphen boundary. ---
$teamRange = Get-ParagraphRangeContaining $d "Abundez-Arce"
$teamRange.Find.Execute("Abundez", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$arceSplitPoint = $teamRange.End
$d.Bookmarks.Add("ZZZTempRunSplit2", $d.Range($arceSplitPoint, $arceSplitPoint))
$d.Bookmarks("ZZZTempRunSplit2").Delete()
